$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking price values in column D are stored as text (inlineStr) in the
# original workbook. Assigning a plain numeric-looking string via .Value would make
# Excel coerce it into a real number (and drop formatting such as trailing zeros),
# so we force text interpretation with a leading apostrophe and then clear the
# resulting cell-level number-format override so the cell keeps its original look.
$priceUpdates = [ordered]@{
    "D2" = "245.70"
    "D3" = "23.84"
    "D4" = "5.325"
    "D5" = "0.05827"
    "D6" = "6.474"
    "D7" = "3.353"
    "D8" = "0.8105"
    "D9" = "0.9210"
    "D10" = "0.1408"
    "D11" = "0.07339"
    "D12" = "0.03080"
    "D13" = "0.03051"
    "D14" = "0.09365"
    "D16" = "0.001570"
    "D17" = "0.04683"
    "D18" = "0.0005993"
    "D19" = "0.006040"
    "D20" = "0.001244"
    "D21" = "0.004690"
    "D22" = "0.00008804"
    "D24" = "2.150"
    "D28" = "0.0002340"
    "D40" = "0.03836"
    "D41" = "0.006343"
    "D43" = "0.003201"
    "D44" = "0.007622"
    "D45" = "0.00005257"
    "D47" = "0.6853"
}

foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.Formula = "'" + $priceUpdates[$cellRef]
    $cell.ClearFormats()
}

# Plain text label updates in column E (these already round-trip as text).
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

Write-Host "Updated $($priceUpdates.Count) price cells and 2 label cells."
